$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Change 1: "Error validation using " -> "Data validation using "
$para5 = $tr.Paragraphs(5, 1)
$run5a = $para5.Runs(1, 1)
$run5a.Text = "Data validation using "

# Change 2: split the "List available container..." paragraph into two
#   paragraph A (new, top level): "Let’s create a process to handle not
#     present libraries and wanted by some app developer."
#   paragraph B (existing, now indented one level): "List available
#     container (templates) for starting a service." as a single run.
$tr2 = $sh.TextFrame.TextRange
$para7 = $tr2.Paragraphs(7, 1)
$para7.InsertBefore("Let’s create a process to handle not present libraries and wanted by some app developer.`r")

$tr3 = $sh.TextFrame.TextRange
$para8 = $tr3.Paragraphs(8, 1)
$para8.IndentLevel = 2

$run8b = $para8.Runs(2, 1)
$run8b.Text = ""

$tr4 = $sh.TextFrame.TextRange
$para8again = $tr4.Paragraphs(8, 1)
$run8a = $para8again.Runs(1, 1)
$run8a.Text = "List available container (templates) for starting a service."
